$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Photometric-Opt")
$ws.Range("D4").Value = 9.919356917703732
$ws.Range("K4").Value = 4.304025912258477
$ws.Range("D5").Value = 9.886936539789014
$ws.Range("K5").Value = 3.981254205101203
$ws.Range("B6").Value = 9.919356917703732
$ws.Range("C6").Value = 9.886936539789014
$ws.Range("E6").Value = 7.117976427563766
$ws.Range("F6").Value = 9.99217094197501
$ws.Range("I6").Value = 4.304025912258477
$ws.Range("J6").Value = 3.981254205101203
$ws.Range("L6").Value = 5.106398705457906
$ws.Range("M6").Value = 3.930611463694983
$ws.Range("D7").Value = 7.117976427563766
$ws.Range("K7").Value = 5.106398705457906
$ws.Range("D8").Value = 9.99217094197501
$ws.Range("K8").Value = 3.930611463694983
$ws.Range("D14").Value = 9.919356917702338
$ws.Range("K14").Value = 4.264097540967208
$ws.Range("D15").Value = 9.886936539793002
$ws.Range("K15").Value = 4.265870609220072
$ws.Range("B16").Value = 9.919356917702338
$ws.Range("C16").Value = 9.886936539793002
$ws.Range("E16").Value = 7.117976427563251
$ws.Range("F16").Value = 9.992170941966545
$ws.Range("I16").Value = 4.264097540967208
$ws.Range("J16").Value = 4.265870609220072
$ws.Range("L16").Value = 4.269155611861901
$ws.Range("M16").Value = 4.26838847448473
$ws.Range("D17").Value = 7.117976427563251
$ws.Range("K17").Value = 4.269155611861901
$ws.Range("D18").Value = 9.992170941966545
$ws.Range("K18").Value = 4.26838847448473

$ws = $wb.Worksheets.Item("Photometric-Pess")
$ws.Range("F6").Value = 10.84918233753177
$ws.Range("M6").Value = 2.617148387586216
$ws.Range("D8").Value = 10.84918233753177
$ws.Range("K8").Value = 2.617148387586216
$ws.Range("F16").Value = 10.84918233753237
$ws.Range("M16").Value = 3.454313084573598
$ws.Range("D18").Value = 10.84918233753237
$ws.Range("K18").Value = 3.454313084573598

$ws = $wb.Worksheets.Item("Spectroscopic-Opt")
$ws.Range("D4").Value = 5.469262880664999
$ws.Range("K4").Value = 4.519503360644053
$ws.Range("D5").Value = 5.497487840891996
$ws.Range("K5").Value = 4.524026768922548
$ws.Range("B6").Value = 5.469262880664999
$ws.Range("C6").Value = 5.497487840891996
$ws.Range("E6").Value = 5.328235491198823
$ws.Range("F6").Value = 5.458721538855309
$ws.Range("I6").Value = 4.519503360644053
$ws.Range("J6").Value = 4.524026768922548
$ws.Range("L6").Value = 2.801845126992966
$ws.Range("M6").Value = 4.501266270591698
$ws.Range("D7").Value = 5.328235491198823
$ws.Range("K7").Value = 2.801845126992966
$ws.Range("D8").Value = 5.458721538855309
$ws.Range("K8").Value = 4.501266270591698
$ws.Range("D14").Value = 5.469262880668691
$ws.Range("K14").Value = 0.3261869951982354
$ws.Range("D15").Value = 5.497487840900265
$ws.Range("K15").Value = 0.3216644176543825
$ws.Range("B16").Value = 5.469262880668691
$ws.Range("C16").Value = 5.497487840900265
$ws.Range("E16").Value = 5.328235491217509
$ws.Range("F16").Value = 5.458721538891282
$ws.Range("I16").Value = 0.3261869951982354
$ws.Range("J16").Value = 0.3216644176543825
$ws.Range("L16").Value = 0.3250320339686323
$ws.Range("M16").Value = 0.3264076035532166
$ws.Range("D17").Value = 5.328235491217509
$ws.Range("K17").Value = 0.3250320339686323
$ws.Range("D18").Value = 5.458721538891282
$ws.Range("K18").Value = 0.3264076035532166
